$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-08-03 21:00:13"

for ($row = 2; $row -le 73; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
